$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.21746643998903
$ws.Range("C2").Value = 10.01825907257359
$ws.Range("D2").Value = 5.978722398905361
$ws.Range("E2").Value = 10.75430003579544
$ws.Range("G2").Value = 37.24697618670292
$ws.Range("H2").Value = 15.51393466838432
$ws.Range("I2").Value = 21.89653471378388
$ws.Range("M2").Value = 15.33594835170326
$ws.Range("B3").Value = 12.6466151517221
$ws.Range("C3").Value = 9.416761198580325
$ws.Range("D3").Value = 5.859837063128313
$ws.Range("E3").Value = 10.65302686721537
$ws.Range("G3").Value = 36.56991370282083
$ws.Range("H3").Value = 15.48809719275254
$ws.Range("I3").Value = 21.87632917135912
$ws.Range("M3").Value = 15.05744823673104
$ws.Range("B4").Value = 12.28571832648902
$ws.Range("C4").Value = 9.028855382972717
$ws.Range("D4").Value = 5.787570046196155
$ws.Range("E4").Value = 10.59410100298315
$ws.Range("G4").Value = 36.16392712441267
$ws.Range("H4").Value = 15.47658731183004
$ws.Range("I4").Value = 21.87066940134487
$ws.Range("M4").Value = 14.88839282652102
$ws.Range("B5").Value = 12.13627143009059
$ws.Range("C5").Value = 8.866227506193647
$ws.Range("D5").Value = 5.758348289760543
$ws.Range("E5").Value = 10.57092810668991
$ws.Range("G5").Value = 36.00117332347321
$ws.Range("H5").Value = 15.4729908959445
$ws.Range("I5").Value = 21.87005285805358
$ws.Range("M5").Value = 14.82008786345369
$ws.Range("B6").Value = 12.111319668254
$ws.Range("C6").Value = 8.838952323199395
$ws.Range("D6").Value = 5.753511086526252
$ws.Range("E6").Value = 10.56713156624905
$ws.Range("G6").Value = 35.97431753604819
$ws.Range("H6").Value = 15.4724597351757
$ws.Range("I6").Value = 21.87005232073612
$ws.Range("M6").Value = 14.80878414416707
$ws.Range("B7").Value = 12.28371213217169
$ws.Range("C7").Value = 9.026680383315792
$ws.Range("D7").Value = 5.787174970948328
$ws.Range("E7").Value = 10.59378505828915
$ws.Range("G7").Value = 36.16172096709215
$ws.Range("H7").Value = 15.47653438219235
$ws.Range("I7").Value = 21.87065425431527
$ws.Range("M7").Value = 14.8874691385731
$ws.Range("B8").Value = 13.02292548283289
$ws.Range("C8").Value = 9.814783533990132
$ws.Range("D8").Value = 5.937605912714576
$ws.Range("E8").Value = 10.71871845255984
$ws.Range("G8").Value = 37.01165165226777
$ws.Range("H8").Value = 15.50412186226225
$ws.Range("I8").Value = 21.88816410163961
$ws.Range("M8").Value = 15.23957630361894
$ws.Range("B9").Value = 14.38120074297948
$ws.Range("C9").Value = 11.26274168591647
$ws.Range("D9").Value = 6.236461624443573
$ws.Range("E9").Value = 10.98853033346244
$ws.Range("G9").Value = 38.74476386768166
$ws.Range("H9").Value = 15.59276872728486
$ws.Range("I9").Value = 21.97625499288157
$ws.Range("M9").Value = 15.94123033716203
$ws.Range("B10").Value = 15.31342246478072
$ws.Range("C10").Value = 12.25807526648683
$ws.Range("D10").Value = 6.455874231815693
$ws.Range("E10").Value = 11.20040879590455
$ws.Range("G10").Value = 40.04401778767318
$ws.Range("H10").Value = 15.6788874765626
$ws.Range("I10").Value = 22.07394614714686
$ws.Range("M10").Value = 16.45794802614346
$ws.Range("B11").Value = 15.72159614537519
$ws.Range("C11").Value = 12.6848147012263
$ws.Range("D11").Value = 6.555152981316099
$ws.Range("E11").Value = 11.29941775484953
$ws.Range("G11").Value = 40.63768761348193
$ws.Range("H11").Value = 15.7225843895598
$ws.Range("I11").Value = 22.12555529899704
$ws.Range("M11").Value = 16.69214333063785
$ws.Range("B12").Value = 15.87376378304813
$ws.Range("C12").Value = 12.84268665528799
$ws.Range("D12").Value = 6.592631110604296
$ws.Range("E12").Value = 11.3372569958838
$ws.Range("G12").Value = 40.86262996542045
$ws.Range("H12").Value = 15.73977626832851
$ws.Range("I12").Value = 22.14612691338613
$ws.Range("M12").Value = 16.78061303212972
$ws.Range("B13").Value = 15.84110006096982
$ws.Range("C13").Value = 12.80885145997627
$ws.Range("D13").Value = 6.584565322457644
$ws.Range("E13").Value = 11.32909271391328
$ws.Range("G13").Value = 40.8141825337313
$ws.Range("H13").Value = 15.73604510797988
$ws.Range("I13").Value = 22.14165077087823
$ws.Range("M13").Value = 16.76157046376919
$ws.Range("B14").Value = 15.73416372937653
$ws.Range("C14").Value = 12.69787746290837
$ws.Range("D14").Value = 6.558238886979097
$ws.Range("E14").Value = 11.30252404161008
$ws.Range("G14").Value = 40.65619247820912
$ws.Range("H14").Value = 15.72398588725367
$ws.Range("I14").Value = 22.12722714916242
$ws.Range("M14").Value = 16.69942660082054
$ws.Range("B15").Value = 15.66834668710279
$ws.Range("C15").Value = 12.62941816866314
$ws.Range("D15").Value = 6.542096855000583
$ws.Range("E15").Value = 11.28629419144042
$ws.Range("G15").Value = 40.55942911079661
$ws.Range("H15").Value = 15.71668306260721
$ws.Range("I15").Value = 22.11852608078687
$ws.Range("M15").Value = 16.66133106548782
$ws.Range("B16").Value = 15.28641739864808
$ws.Range("C16").Value = 12.229664497722
$ws.Range("D16").Value = 6.449371943729855
$ws.Range("E16").Value = 11.19398852892282
$ws.Range("G16").Value = 40.00525180422031
$ws.Range("H16").Value = 15.67612232649222
$ws.Range("I16").Value = 22.0707174480765
$ws.Range("M16").Value = 16.44261758904188
$ws.Range("B17").Value = 15.0479582796661
$ws.Range("C17").Value = 11.97777258383334
$ws.Range("D17").Value = 6.392322600512798
$ws.Range("E17").Value = 11.13801254959915
$ws.Range("G17").Value = 39.66577840932168
$ws.Range("H17").Value = 15.65239414876227
$ws.Range("I17").Value = 22.04322328783405
$ws.Range("M17").Value = 16.30815707845287
$ws.Range("B18").Value = 14.90931197809293
$ws.Range("C18").Value = 11.83043856010316
$ws.Range("D18").Value = 6.359461404449725
$ws.Range("E18").Value = 11.10606555417
$ws.Range("G18").Value = 39.47078431474596
$ws.Range("H18").Value = 15.63917226166757
$ws.Range("I18").Value = 22.0280844133565
$ws.Range("M18").Value = 16.23074299203414
$ws.Range("B19").Value = 14.86211637935203
$ws.Range("C19").Value = 11.78013197111425
$ws.Range("D19").Value = 6.348328187905061
$ws.Range("E19").Value = 11.09529252436521
$ws.Range("G19").Value = 39.4048153103583
$ws.Range("H19").Value = 15.6347688561601
$ws.Range("I19").Value = 22.02307463799475
$ws.Range("M19").Value = 16.20452170756945
$ws.Range("B20").Value = 15.07349783674054
$ws.Range("C20").Value = 12.00484058343153
$ws.Range("D20").Value = 6.398400856166632
$ws.Range("E20").Value = 11.14394574740626
$ws.Range("G20").Value = 39.70189062458924
$ws.Range("H20").Value = 15.65487600635849
$ws.Range("I20").Value = 22.04608023701674
$ws.Range("M20").Value = 16.32247916313528
$ws.Range("B21").Value = 15.76563943902843
$ws.Range("C21").Value = 12.73057415333808
$ws.Range("D21").Value = 6.565975070623284
$ws.Range("E21").Value = 11.31031874304886
$ws.Range("G21").Value = 40.70259630776129
$ws.Range("H21").Value = 15.72751052064603
$ws.Range("I21").Value = 22.13143583621994
$ws.Range("M21").Value = 16.71768629833565
$ws.Range("B22").Value = 16.2039710296171
$ws.Range("C22").Value = 13.18317996311126
$ws.Range("D22").Value = 6.674798024725069
$ws.Range("E22").Value = 11.42105931561199
$ws.Range("G22").Value = 41.35728012449592
$ws.Range("H22").Value = 15.77873656567354
$ws.Range("I22").Value = 22.1932122620515
$ws.Range("M22").Value = 16.97468482064063
$ws.Range("B23").Value = 15.97134115181522
$ws.Range("C23").Value = 12.94359500018997
$ws.Range("D23").Value = 6.61679352055878
$ws.Range("E23").Value = 11.36178170480809
$ws.Range("G23").Value = 41.00788108188259
$ws.Range("H23").Value = 15.75105469907747
$ws.Range("I23").Value = 22.15969407343545
$ws.Range("M23").Value = 16.83766705556177
$ws.Range("B24").Value = 15.06195623143658
$ws.Range("C24").Value = 11.99261098535224
$ws.Range("D24").Value = 6.395653069244453
$ws.Range("E24").Value = 11.14126261643876
$ws.Range("G24").Value = 39.68556373231133
$ws.Range("H24").Value = 15.65375265069545
$ws.Range("I24").Value = 22.0447865298081
$ws.Range("M24").Value = 16.31600448727747
$ws.Range("B25").Value = 14.02465810748618
$ws.Range("C25").Value = 10.87297481648394
$ws.Range("D25").Value = 6.155464041716987
$ws.Range("E25").Value = 10.91302933865821
$ws.Range("G25").Value = 38.27031108261369
$ws.Range("H25").Value = 15.56509170988685
$ws.Range("I25").Value = 21.94664187672292
$ws.Range("M25").Value = 15.75080490952943
